# Add a "Role" column to the exported claims sheet, right after "MiddleName"
# and before "Branch" (i.e. a new column F), shifting the remaining columns
# (Branch..Month of Claim) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; everything from old column F onward
# (Branch, Sol ID, Weekday, Weekend, Shift, Amount, Status, Month of Claim)
# shifts right by one, carrying its formatting with it.
$ws.Range("F1").EntireColumn.Insert()

# Header for the new column
$ws.Range("F1").Value = "Role"

# Values for the new column, one per data row
$ws.Range("F2").Value = "Service Executive (Financial)"
$ws.Range("F3").Value = "Service Executive (Financial)"
$ws.Range("F4").Value = "Service Executive (Non-Financial)"
